# Landscaping Data.xlsx update
# - Revise several "Growth" (column H) readings for existing rows
# - Append 7 new observation rows (611-617) for 8/5 (serial 45874)
# - Move the active selection to I618 (just past the new data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Corrected Growth (column H) values for existing rows.
#    A few of these cells previously held formulas (=1/3, =2/3, =8/3); the
#    corrected entries are plain numbers, so assigning .Value replaces the
#    formula outright.
# ---------------------------------------------------------------------------
$hUpdates = @{
    43  = 0.25
    51  = 0.3
    52  = 0.6
    53  = 0.55000000000000004
    54  = 0.8
    55  = 0.9
    56  = 1.1000000000000001
    57  = 4
    116 = 0
    162 = 0.5
    178 = 0
    268 = 0.5
    271 = 0.75
    272 = 0.8
    273 = 2.25
    274 = 5.5
    275 = 0.75
    281 = 7.75
}

foreach ($row in $hUpdates.Keys) {
    $ws.Cells.Item($row, 8).Value = $hUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2. Append the new rows of collected data (rows 611-617).
#    First clone the formatting of the last existing row (610) down across
#    the new block so date styling etc. carries through, then fill in values.
# ---------------------------------------------------------------------------
$lastRow = 610
$newLastRow = 617

$ws.Range("A610:T610").Copy() | Out-Null
$ws.Range("A611:T617").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row=611; Date=45874; PlantType="Flowering";    PlantSize="Large";  Low=63; High=85; Rain=0; Growth=0.1;  Pruned="No"; Quadrant=2; Shade="Bright";  UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
    @{ Row=612; Date=45874; PlantType="Nonflowering"; PlantSize="Medium"; Low=63; High=85; Rain=0; Growth=0;    Pruned="No"; Quadrant=3; Shade="Bright";  UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
    @{ Row=613; Date=45874; PlantType="Nonflowering"; PlantSize="Small";  Low=63; High=85; Rain=0; Growth=0;    Pruned="No"; Quadrant=3; Shade="Neutral"; UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
    @{ Row=614; Date=45874; PlantType="Nonflowering"; PlantSize="Medium"; Low=63; High=85; Rain=0; Growth=0;    Pruned="No"; Quadrant=3; Shade="Dark";    UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
    @{ Row=615; Date=45874; PlantType="Nonflowering"; PlantSize="Medium"; Low=63; High=85; Rain=0; Growth=0.2;  Pruned="No"; Quadrant=3; Shade="Neutral"; UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
    @{ Row=616; Date=45874; PlantType="Nonflowering"; PlantSize="Large";  Low=63; High=85; Rain=0; Growth=0.5;  Pruned="No"; Quadrant=4; Shade="Dark";    UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
    @{ Row=617; Date=45874; PlantType="Tree";          PlantSize="Medium"; Low=63; High=85; Rain=0; Growth=0.75; Pruned="No"; Quadrant=1; Shade="Bright";  UV=8; Humidity=0.56000000000000005; DewPoint=61; Pressure=30.27; WindGust=9; CloudCover=0.35; Visibility=9.9; AQI=53; Pollen=6 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.PlantType
    $ws.Cells.Item($row, 3).Value = $r.PlantSize
    $ws.Cells.Item($row, 4).Value = $r.Low
    $ws.Cells.Item($row, 5).Value = $r.High
    $ws.Range("F$row").Formula = "=ABS(D$row-E$row)"
    $ws.Cells.Item($row, 7).Value  = $r.Rain
    $ws.Cells.Item($row, 8).Value  = $r.Growth
    $ws.Cells.Item($row, 9).Value  = $r.Pruned
    $ws.Cells.Item($row, 10).Value = $r.Quadrant
    $ws.Cells.Item($row, 11).Value = $r.Shade
    $ws.Cells.Item($row, 12).Value = $r.UV
    $ws.Cells.Item($row, 13).Value = $r.Humidity
    $ws.Cells.Item($row, 14).Value = $r.DewPoint
    $ws.Cells.Item($row, 15).Value = $r.Pressure
    $ws.Cells.Item($row, 16).Value = $r.WindGust
    $ws.Cells.Item($row, 17).Value = $r.CloudCover
    $ws.Cells.Item($row, 18).Value = $r.Visibility
    $ws.Cells.Item($row, 19).Value = $r.AQI
    $ws.Cells.Item($row, 20).Value = $r.Pollen
}

# ---------------------------------------------------------------------------
# 3. Move the selection to reflect where entry would continue next (I618).
# ---------------------------------------------------------------------------
$ws.Range("I618").Select() | Out-Null
